# Edit student feedback commit numbers.
#
# Two numeric changes are made to the document's body text:
#   1. "Changes made at commit 301 for the final draft."  -> commit 303
#      (end of the "Part 2: Data Collection and Cleaning Blog" section)
#   2. "Final draft is commit 25."                         -> commit 26
#      (end of the "Part 4: Streamlit Dashboard" section)
#
# Both substrings are unique within the document, so a simple Find/Replace
# on the whole document content is sufficient and safe.

$d = $word.ActiveDocument

# wdReplaceAll = 2
$wdReplaceAll = 2

# --- Change 1: commit 301 -> commit 303 -----------------------------------
$range1 = $d.Content
$range1.Find.Execute(
    "commit 301 for the final draft.",  # FindText
    $true,                              # MatchCase
    $true,                              # MatchWholeWord
    $false,                             # MatchWildcards
    $false,                             # MatchSoundsLike
    $false,                             # MatchAllWordForms
    $true,                              # Forward
    1,                                  # Wrap (wdFindContinue)
    $false,                             # Format
    "commit 303 for the final draft.",  # ReplaceWith
    $wdReplaceAll
) | Out-Null

# --- Change 2: commit 25 -> commit 26 (final draft commit number) ---------
$range2 = $d.Content
$range2.Find.Execute(
    "Final draft is commit 25.",
    $true,
    $true,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "Final draft is commit 26.",
    $wdReplaceAll
) | Out-Null
